# Generate Report for handoff
# The file "7593a10e-691d-4a59-b935-bd2d1ef4e50c.md" is now ready for handoff
# (new handoff generated), so update its Status to "Ready for handoff" on all
# sheets, and stamp the new "Latest Handoff Datetime" on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 7593a10e...md is row 3 ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row for 7593a10e...md is row 3 ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("D3").Value = "2016-01-13 15:48:43"

# --- de-de sheet: row for 7593a10e...md is row 3 ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("B3").Value = "Ready for handoff"
$de.Range("D3").Value = "2016-01-13 15:48:54"
